$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.02328126719340038
$ws.Range("J2").Value = 0.02328126719340038
$ws.Range("M2").Value = 2.815739333333333
$ws.Range("N2").Value = 8.447217999999999
$ws.Range("O2").Value = 0.07700398964630729
$ws.Range("P2").Value = 0.07700398964630729
$ws.Range("Q2").Value = 0.05937549532199999
$ws.Range("R2").Value = 0.5343794578979999
$ws.Range("S2").Value = 0.001792750457913517
$ws.Range("T2").Value = 0.001792750457913517

# Row 3
$ws.Range("I3").Value = 0.02328126719340038
$ws.Range("J3").Value = 0.02328126719340038
$ws.Range("O3").Value = 0.1324338085883186
$ws.Range("P3").Value = 0.1324338085883186
$ws.Range("S3").Value = 0.003083226883184287
$ws.Range("T3").Value = 0.003083226883184287

# Row 4
$ws.Range("I4").Value = 0.02328126719340038
$ws.Range("J4").Value = 0.02328126719340038
$ws.Range("M4").Value = 5.537790999999999
$ws.Range("N4").Value = 16.613373
$ws.Range("O4").Value = 0.1514458372546134
$ws.Range("P4").Value = 0.1514458372546134
$ws.Range("Q4").Value = 0.116775398817
$ws.Range("R4").Value = 1.050978589353
$ws.Range("S4").Value = 0.003525851002452885
$ws.Range("T4").Value = 0.003525851002452886

# Row 5
$ws.Range("I5").Value = 0.02328126719340038
$ws.Range("J5").Value = 0.02328126719340038
$ws.Range("M5").Value = 1.188595666666667
$ws.Range("N5").Value = 3.565787
$ws.Range("O5").Value = 0.03250535563648733
$ws.Range("P5").Value = 0.03250535563648733
$ws.Range("Q5").Value = 0.025063916823
$ws.Range("R5").Value = 0.225575251407
$ws.Range("S5").Value = 0.0007567658697895645
$ws.Range("T5").Value = 0.0007567658697895647

# Row 6
$ws.Range("I6").Value = 0.02328126719340038
$ws.Range("J6").Value = 0.02328126719340038
$ws.Range("M6").Value = 18.85109966666667
$ws.Range("N6").Value = 56.553299
$ws.Range("O6").Value = 0.5155341854158992
$ws.Range("P6").Value = 0.5155341854158992
$ws.Range("Q6").Value = 0.397513138671
$ws.Range("R6").Value = 3.577618248039
$ws.Range("S6").Value = 0.01200228911799956
$ws.Range("T6").Value = 0.01200228911799957

# Row 7
$ws.Range("I7").Value = 0.02328126719340038
$ws.Range("J7").Value = 0.02328126719340038
$ws.Range("M7").Value = 3.330328666666666
$ws.Range("N7").Value = 9.990985999999999
$ws.Range("O7").Value = 0.09107682345837424
$ws.Range("P7").Value = 0.09107682345837424
$ws.Range("Q7").Value = 0.07022664059399999
$ws.Range("R7").Value = 0.6320397653459999
$ws.Range("S7").Value = 0.002120383862060567
$ws.Range("T7").Value = 0.002120383862060567

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8846626666666667
$ws.Range("H8").Value = 2.653988
$ws.Range("I8").Value = 0.9767187328065996
$ws.Range("J8").Value = 0.9767187328065997
$ws.Range("M8").Value = 2.815739333333333
$ws.Range("N8").Value = 8.447217999999999
$ws.Range("O8").Value = 0.07700398964630729
$ws.Range("P8").Value = 0.07700398964630729
$ws.Range("Q8").Value = 2.490979467264889
$ws.Range("R8").Value = 22.418815205384
$ws.Range("S8").Value = 0.07521123918839377
$ws.Range("T8").Value = 0.07521123918839379

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8846626666666667
$ws.Range("H9").Value = 2.653988
$ws.Range("I9").Value = 0.9767187328065996
$ws.Range("J9").Value = 0.9767187328065997
$ws.Range("O9").Value = 0.1324338085883186
$ws.Range("P9").Value = 0.1324338085883186
$ws.Range("Q9").Value = 4.284062416511556
$ws.Range("R9").Value = 38.556561748604
$ws.Range("S9").Value = 0.1293505817051343
$ws.Range("T9").Value = 0.1293505817051343

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8846626666666667
$ws.Range("H10").Value = 2.653988
$ws.Range("I10").Value = 0.9767187328065996
$ws.Range("J10").Value = 0.9767187328065997
$ws.Range("M10").Value = 5.537790999999999
$ws.Range("N10").Value = 16.613373
$ws.Range("O10").Value = 0.1514458372546134
$ws.Range("P10").Value = 0.1514458372546134
$ws.Range("Q10").Value = 4.899076953502666
$ws.Range("R10").Value = 44.091692581524
$ws.Range("S10").Value = 0.1479199862521605
$ws.Range("T10").Value = 0.1479199862521606

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.8846626666666667
$ws.Range("H11").Value = 2.653988
$ws.Range("I11").Value = 0.9767187328065996
$ws.Range("J11").Value = 0.9767187328065997
$ws.Range("M11").Value = 1.188595666666667
$ws.Range("N11").Value = 3.565787
$ws.Range("O11").Value = 0.03250535563648733
$ws.Range("P11").Value = 0.03250535563648733
$ws.Range("Q11").Value = 1.051506212061778
$ws.Range("R11").Value = 9.463555908556001
$ws.Range("S11").Value = 0.03174858976669776
$ws.Range("T11").Value = 0.03174858976669777

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.8846626666666667
$ws.Range("H12").Value = 2.653988
$ws.Range("I12").Value = 0.9767187328065996
$ws.Range("J12").Value = 0.9767187328065997
$ws.Range("M12").Value = 18.85109966666667
$ws.Range("N12").Value = 56.553299
$ws.Range("O12").Value = 0.5155341854158992
$ws.Range("P12").Value = 0.5155341854158992
$ws.Range("Q12").Value = 16.67686410071245
$ws.Range("R12").Value = 150.091776906412
$ws.Range("S12").Value = 0.5035318962978996
$ws.Range("T12").Value = 0.5035318962978997

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.8846626666666667
$ws.Range("H13").Value = 2.653988
$ws.Range("I13").Value = 0.9767187328065996
$ws.Range("J13").Value = 0.9767187328065997
$ws.Range("M13").Value = 3.330328666666666
$ws.Range("N13").Value = 9.990985999999999
$ws.Range("O13").Value = 0.09107682345837424
$ws.Range("P13").Value = 0.09107682345837424
$ws.Range("Q13").Value = 2.946217439129778
$ws.Range("R13").Value = 26.515956952168
$ws.Range("S13").Value = 0.08895643959631368
$ws.Range("T13").Value = 0.08895643959631369
